# Update NATMI ligand-receptor (Ntn3-Neo1) output with re-computed TPM-based values.
# The workbook stores pre-computed numeric results (no formulas), so the cells are
# overwritten directly with the values produced by the updated TPM pipeline.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.4940636666666666
$ws.Range("H2").Value = 1.482191
$ws.Range("I2").Value = 0.1416906061387336
$ws.Range("J2").Value = 0.1416906061387335
$ws.Range("M2").Value = 2.733663333333333
$ws.Range("N2").Value = 8.200989999999999
$ws.Range("O2").Value = 0.04037266183309663
$ws.Range("P2").Value = 0.04037266183309663
$ws.Range("Q2").Value = 1.350603729898889
$ws.Range("R2").Value = 12.15543356909
$ws.Range("S2").Value = 0.005720426926565576
$ws.Range("T2").Value = 0.005720426926565575
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.4940636666666666
$ws.Range("H3").Value = 1.482191
$ws.Range("I3").Value = 0.1416906061387336
$ws.Range("J3").Value = 0.1416906061387335
$ws.Range("O3").Value = 0.6389522306252696
$ws.Range("P3").Value = 0.6389522306252696
$ws.Range("Q3").Value = 21.37513918396778
$ws.Range("R3").Value = 192.37625265571
$ws.Range("S3").Value = 0.09053352885099034
$ws.Range("T3").Value = 0.09053352885099032
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.4940636666666666
$ws.Range("H4").Value = 1.482191
$ws.Range("I4").Value = 0.1416906061387336
$ws.Range("J4").Value = 0.1416906061387335
$ws.Range("M4").Value = 21.46453166666667
$ws.Range("N4").Value = 64.393595
$ws.Range("O4").Value = 0.3170032929137071
$ws.Range("P4").Value = 0.317003292913707
$ws.Range("Q4").Value = 10.60484521851611
$ws.Range("R4").Value = 95.44360696664499
$ws.Range("S4").Value = 0.04491638872091767
$ws.Range("T4").Value = 0.04491638872091765
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.4940636666666666
$ws.Range("H5").Value = 1.482191
$ws.Range("I5").Value = 0.1416906061387336
$ws.Range("J5").Value = 0.1416906061387335
$ws.Range("M5").Value = 0.2486213333333333
$ws.Range("N5").Value = 0.745864
$ws.Range("O5").Value = 0.003671814627926724
$ws.Range("P5").Value = 0.003671814627926724
$ws.Range("Q5").Value = 0.1228347675582222
$ws.Range("R5").Value = 1.105512908024
$ws.Range("S5").Value = 0.0005202616402600061
$ws.Range("T5").Value = 0.000520261640260006
$ws.Range("H6").Value = 6.480663
$ws.Range("I6").Value = 0.6195214170446748
$ws.Range("J6").Value = 0.6195214170446747
$ws.Range("M6").Value = 2.733663333333333
$ws.Range("N6").Value = 8.200989999999999
$ws.Range("O6").Value = 0.04037266183309663
$ws.Range("P6").Value = 0.04037266183309663
$ws.Range("Q6").Value = 5.905316939596665
$ws.Range("R6").Value = 53.14785245637
$ws.Range("S6").Value = 0.02501172866870548
$ws.Range("T6").Value = 0.02501172866870548
$ws.Range("H7").Value = 6.480663
$ws.Range("I7").Value = 0.6195214170446748
$ws.Range("J7").Value = 0.6195214170446747
$ws.Range("O7").Value = 0.6389522306252696
$ws.Range("P7").Value = 0.6389522306252696
$ws.Range("Q7").Value = 93.45966453000332
$ws.Range("R7").Value = 841.1369807700299
$ws.Range("S7").Value = 0.3958445913408229
$ws.Range("T7").Value = 0.3958445913408228
$ws.Range("H8").Value = 6.480663
$ws.Range("I8").Value = 0.6195214170446748
$ws.Range("J8").Value = 0.6195214170446747
$ws.Range("M8").Value = 21.46453166666667
$ws.Range("N8").Value = 64.393595
$ws.Range("O8").Value = 0.3170032929137071
$ws.Range("P8").Value = 0.317003292913707
$ws.Range("Q8").Value = 46.36813206149834
$ws.Range("R8").Value = 417.313188553485
$ws.Range("S8").Value = 0.1963903292337279
$ws.Range("T8").Value = 0.1963903292337279
$ws.Range("H9").Value = 6.480663
$ws.Range("I9").Value = 0.6195214170446748
$ws.Range("J9").Value = 0.6195214170446747
$ws.Range("M9").Value = 0.2486213333333333
$ws.Range("N9").Value = 0.745864
$ws.Range("O9").Value = 0.003671814627926724
$ws.Range("P9").Value = 0.003671814627926724
$ws.Range("Q9").Value = 0.5370770253146666
$ws.Range("R9").Value = 4.833693227832
$ws.Range("S9").Value = 0.00227476780141853
$ws.Range("T9").Value = 0.002274767801418529
$ws.Range("G10").Value = 0.7459539999999999
$ws.Range("H10").Value = 2.237862
$ws.Range("I10").Value = 0.2139292596128559
$ws.Range("J10").Value = 0.2139292596128559
$ws.Range("M10").Value = 2.733663333333333
$ws.Range("N10").Value = 8.200989999999999
$ws.Range("O10").Value = 0.04037266183309663
$ws.Range("P10").Value = 0.04037266183309663
$ws.Range("Q10").Value = 2.039187098153333
$ws.Range("R10").Value = 18.35268388338
$ws.Range("S10").Value = 0.00863689365455457
$ws.Range("T10").Value = 0.00863689365455457
$ws.Range("G11").Value = 0.7459539999999999
$ws.Range("H11").Value = 2.237862
$ws.Range("I11").Value = 0.2139292596128559
$ws.Range("J11").Value = 0.2139292596128559
$ws.Range("O11").Value = 0.6389522306252696
$ws.Range("P11").Value = 0.6389522306252696
$ws.Range("Q11").Value = 32.27290661224666
$ws.Range("R11").Value = 290.4561595102199
$ws.Range("S11").Value = 0.1366905776256467
$ws.Range("T11").Value = 0.1366905776256467
$ws.Range("G12").Value = 0.7459539999999999
$ws.Range("H12").Value = 2.237862
$ws.Range("I12").Value = 0.2139292596128559
$ws.Range("J12").Value = 0.2139292596128559
$ws.Range("M12").Value = 21.46453166666667
$ws.Range("N12").Value = 64.393595
$ws.Range("O12").Value = 0.3170032929137071
$ws.Range("P12").Value = 0.317003292913707
$ws.Range("Q12").Value = 16.01155325487667
$ws.Range("R12").Value = 144.10397929389
$ws.Range("S12").Value = 0.06781627974786666
$ws.Range("T12").Value = 0.06781627974786665
$ws.Range("G13").Value = 0.7459539999999999
$ws.Range("H13").Value = 2.237862
$ws.Range("I13").Value = 0.2139292596128559
$ws.Range("J13").Value = 0.2139292596128559
$ws.Range("M13").Value = 0.2486213333333333
$ws.Range("N13").Value = 0.745864
$ws.Range("O13").Value = 0.003671814627926724
$ws.Range("P13").Value = 0.003671814627926724
$ws.Range("Q13").Value = 0.1854600780853333
$ws.Range("R13").Value = 1.669140702768
$ws.Range("S13").Value = 0.0007855085847880182
$ws.Range("T13").Value = 0.0007855085847880182
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.08668033333333335
$ws.Range("H14").Value = 0.260041
$ws.Range("I14").Value = 0.02485871720373584
$ws.Range("J14").Value = 0.02485871720373583
$ws.Range("M14").Value = 2.733663333333333
$ws.Range("N14").Value = 8.200989999999999
$ws.Range("O14").Value = 0.04037266183309663
$ws.Range("P14").Value = 0.04037266183309663
$ws.Range("Q14").Value = 0.2369548489544444
$ws.Range("R14").Value = 2.13259364059
$ws.Range("S14").Value = 0.001003612583271008
$ws.Range("T14").Value = 0.001003612583271008
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.08668033333333335
$ws.Range("H15").Value = 0.260041
$ws.Range("I15").Value = 0.02485871720373584
$ws.Range("J15").Value = 0.02485871720373583
$ws.Range("O15").Value = 0.6389522306252696
$ws.Range("P15").Value = 0.6389522306252696
$ws.Range("Q15").Value = 3.750132451578889
$ws.Range("R15").Value = 33.75119206421
$ws.Range("S15").Value = 0.01588353280780978
$ws.Range("T15").Value = 0.01588353280780978
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.08668033333333335
$ws.Range("H16").Value = 0.260041
$ws.Range("I16").Value = 0.02485871720373584
$ws.Range("J16").Value = 0.02485871720373583
$ws.Range("M16").Value = 21.46453166666667
$ws.Range("N16").Value = 64.393595
$ws.Range("O16").Value = 0.3170032929137071
$ws.Range("P16").Value = 0.317003292913707
$ws.Range("Q16").Value = 1.860552759710556
$ws.Range("R16").Value = 16.744974837395
$ws.Range("S16").Value = 0.007880295211194882
$ws.Range("T16").Value = 0.007880295211194878
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.08668033333333335
$ws.Range("H17").Value = 0.260041
$ws.Range("I17").Value = 0.02485871720373584
$ws.Range("J17").Value = 0.02485871720373583
$ws.Range("M17").Value = 0.2486213333333333
$ws.Range("N17").Value = 0.745864
$ws.Range("O17").Value = 0.003671814627926724
$ws.Range("P17").Value = 0.003671814627926724
$ws.Range("Q17").Value = 0.02155058004711111
$ws.Range("R17").Value = 0.193955220424
$ws.Range("S17").Value = 0.00009127660146017096
$ws.Range("T17").Value = 0.00009127660146017095
